$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 32675.883
$ws.Cells.Item(64, 9).Value = 3382.5
$ws.Cells.Item(64, 10).Value = 102980
$ws.Cells.Item(64, 11).Value = 3382.5
$ws.Cells.Item(64, 12).Value = 102980
$ws.Cells.Item(64, 13).Value = -3134.5
$ws.Cells.Item(64, 14).Value = -103476

$ws.Cells.Item(67, 8).Value = 32675.883
$ws.Cells.Item(67, 9).Value = 3382.5
$ws.Cells.Item(67, 10).Value = 102980
$ws.Cells.Item(67, 11).Value = 3382.5
$ws.Cells.Item(67, 12).Value = 102980
$ws.Cells.Item(67, 13).Value = -2524.5
$ws.Cells.Item(67, 14).Value = -104696

$ws.Cells.Item(74, 8).Value = 2276039.8
$ws.Cells.Item(74, 9).Value = 2276039.8
$ws.Cells.Item(74, 11).Value = 2276039.8
$ws.Cells.Item(74, 13).Value = -2275103.8

$ws.Cells.Item(76, 8).Value = 27030162
$ws.Cells.Item(76, 9).Value = 31252940
$ws.Cells.Item(76, 10).Value = 4388.8
$ws.Cells.Item(76, 11).Value = 31252940
$ws.Cells.Item(76, 12).Value = 4388.8
$ws.Cells.Item(76, 13).Value = -31252625
$ws.Cells.Item(76, 14).Value = -5018.8

$ws.Cells.Item(77, 8).Value = 2276039.8
$ws.Cells.Item(77, 9).Value = 2276039.8
$ws.Cells.Item(77, 11).Value = 11380199
$ws.Cells.Item(77, 13).Value = -11375519

$ws.Cells.Item(79, 8).Value = 27030162
$ws.Cells.Item(79, 9).Value = 31252940
$ws.Cells.Item(79, 10).Value = 4388.8
$ws.Cells.Item(79, 11).Value = 31252940
$ws.Cells.Item(79, 12).Value = 4388.8
$ws.Cells.Item(79, 13).Value = -31251848
$ws.Cells.Item(79, 14).Value = -6572.8

$ws.Cells.Item(92, 8).Value = 71428770
$ws.Cells.Item(92, 9).Value = 100000180
$ws.Cells.Item(92, 10).Value = 250
$ws.Cells.Item(92, 11).Value = 100000180
$ws.Cells.Item(92, 12).Value = 250
$ws.Cells.Item(92, 13).Value = -99998932
$ws.Cells.Item(92, 14).Value = -2746

$ws.Cells.Item(97, 8).Value = 9703.333000000001
$ws.Cells.Item(97, 10).Value = 9703.333000000001
$ws.Cells.Item(97, 12).Value = 29109.999
$ws.Cells.Item(97, 14).Value = -30101.999

$ws.Cells.Item(99, 8).Value = 104169160
$ws.Cells.Item(99, 9).Value = 50004280
$ws.Cells.Item(99, 10).Value = 142858370
$ws.Cells.Item(99, 11).Value = 150012840
$ws.Cells.Item(99, 12).Value = 428575110
$ws.Cells.Item(99, 13).Value = -150011342
$ws.Cells.Item(99, 14).Value = -428578106

$ws.Cells.Item(100, 8).Value = 5082.8823
$ws.Cells.Item(100, 9).Value = 3808.5833
$ws.Cells.Item(100, 10).Value = 8141.2
$ws.Cells.Item(100, 11).Value = 3808.5833
$ws.Cells.Item(100, 12).Value = 8141.2
$ws.Cells.Item(100, 13).Value = -3267.5833
$ws.Cells.Item(100, 14).Value = -9223.200000000001

$ws.Cells.Item(112, 8).Value = 1295.09
$ws.Cells.Item(112, 10).Value = 1332.3549
$ws.Cells.Item(112, 12).Value = 3997.0647
$ws.Cells.Item(112, 14).Value = -6213.0647

$ws.Cells.Item(113, 8).Value = 1855.68
$ws.Cells.Item(113, 9).Value = 1769.7646
$ws.Cells.Item(113, 10).Value = 2038.25
$ws.Cells.Item(113, 11).Value = 1769.7646
$ws.Cells.Item(113, 12).Value = 2038.25
$ws.Cells.Item(113, 13).Value = 1484.2354
$ws.Cells.Item(113, 14).Value = -8546.25

$ws.Cells.Item(132, 8).Value = 239258.5
$ws.Cells.Item(132, 9).Value = 1112.5834
$ws.Cells.Item(132, 10).Value = 1668134
$ws.Cells.Item(132, 11).Value = 3337.7502
$ws.Cells.Item(132, 12).Value = 5004402
$ws.Cells.Item(132, 13).Value = -807.7501999999999
$ws.Cells.Item(132, 14).Value = -5009462

$ws.Cells.Item(138, 8).Value = 1718.9734
$ws.Cells.Item(138, 9).Value = 863.8
$ws.Cells.Item(138, 10).Value = 2029.9454
$ws.Cells.Item(138, 11).Value = 2591.4
$ws.Cells.Item(138, 12).Value = 6089.8362
$ws.Cells.Item(138, 13).Value = 2548.6
$ws.Cells.Item(138, 14).Value = -16369.8362

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 829.625
$ws.Cells.Item(2, 9).Value = 801
$ws.Cells.Item(2, 10).Value = 892.6
$ws.Cells.Item(2, 11).Value = 801
$ws.Cells.Item(2, 12).Value = 892.6
$ws.Cells.Item(2, 13).Value = -688
$ws.Cells.Item(2, 14).Value = -1118.6

$ws.Cells.Item(45, 8).Value = 1918.0476
$ws.Cells.Item(45, 9).Value = 1463.4706
$ws.Cells.Item(45, 11).Value = 1463.4706
$ws.Cells.Item(45, 13).Value = -1086.4706

$ws.Cells.Item(97, 8).Value = 1225.8387
$ws.Cells.Item(97, 9).Value = 865.3461
$ws.Cells.Item(97, 10).Value = 3100.4
$ws.Cells.Item(97, 11).Value = 865.3461
$ws.Cells.Item(97, 12).Value = 3100.4
$ws.Cells.Item(97, 13).Value = -369.3461
$ws.Cells.Item(97, 14).Value = -4092.4

$ws.Cells.Item(116, 8).Value = 829.625
$ws.Cells.Item(116, 9).Value = 801
$ws.Cells.Item(116, 10).Value = 892.6
$ws.Cells.Item(116, 11).Value = 801
$ws.Cells.Item(116, 12).Value = 892.6
$ws.Cells.Item(116, 13).Value = 1493
$ws.Cells.Item(116, 14).Value = -5480.6

$ws.Cells.Item(122, 8).Value = 1390.5834
$ws.Cells.Item(122, 9).Value = 1065.2
$ws.Cells.Item(122, 11).Value = 3195.6
$ws.Cells.Item(122, 13).Value = -745.6000000000004

$ws.Cells.Item(133, 8).Value = 95000
$ws.Cells.Item(133, 10).Value = 95000
$ws.Cells.Item(133, 12).Value = 95000
$ws.Cells.Item(133, 14).Value = -100060

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 829.625
$ws.Cells.Item(3, 9).Value = 801
$ws.Cells.Item(3, 10).Value = 892.6
$ws.Cells.Item(3, 11).Value = 801
$ws.Cells.Item(3, 12).Value = 892.6
$ws.Cells.Item(3, 13).Value = -687
$ws.Cells.Item(3, 14).Value = -1120.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 90912610
$ws.Cells.Item(86, 9).Value = 200004180
$ws.Cells.Item(86, 10).Value = 2966.6667
$ws.Cells.Item(86, 11).Value = 200004180
$ws.Cells.Item(86, 12).Value = 2966.6667
$ws.Cells.Item(86, 13).Value = -200003057
$ws.Cells.Item(86, 14).Value = -5212.6667

$ws.Cells.Item(89, 8).Value = 90912610
$ws.Cells.Item(89, 9).Value = 200004180
$ws.Cells.Item(89, 10).Value = 2966.6667
$ws.Cells.Item(89, 11).Value = 1000020900
$ws.Cells.Item(89, 12).Value = 14833.3335
$ws.Cells.Item(89, 13).Value = -1000015284
$ws.Cells.Item(89, 14).Value = -26065.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 5028.087
$ws.Cells.Item(5, 9).Value = 444.42856
$ws.Cells.Item(5, 10).Value = 7033.4375
$ws.Cells.Item(5, 11).Value = 1333.28568
$ws.Cells.Item(5, 12).Value = 21100.3125
$ws.Cells.Item(5, 13).Value = -1221.28568
$ws.Cells.Item(5, 14).Value = -21324.3125

$ws.Cells.Item(122, 8).Value = 368.83334
$ws.Cells.Item(122, 9).Value = 322.6
$ws.Cells.Item(122, 11).Value = 2903.4
$ws.Cells.Item(122, 13).Value = -453.4000000000001

$ws.Cells.Item(134, 8).Value = 4809.15
$ws.Cells.Item(134, 9).Value = 1839.1111
$ws.Cells.Item(134, 10).Value = 7239.1816
$ws.Cells.Item(134, 11).Value = 5517.3333
$ws.Cells.Item(134, 12).Value = 21717.5448
$ws.Cells.Item(134, 13).Value = -447.3333000000002
$ws.Cells.Item(134, 14).Value = -31857.5448

$ws.Cells.Item(135, 8).Value = 5028.087
$ws.Cells.Item(135, 9).Value = 444.42856
$ws.Cells.Item(135, 10).Value = 7033.4375
$ws.Cells.Item(135, 11).Value = 3999.85704
$ws.Cells.Item(135, 12).Value = 63300.9375
$ws.Cells.Item(135, 13).Value = -1464.85704
$ws.Cells.Item(135, 14).Value = -68370.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 53344.332
$ws.Cells.Item(52, 10).Value = 53344.332
$ws.Cells.Item(52, 12).Value = 53344.332
$ws.Cells.Item(52, 14).Value = -53862.332

$ws.Cells.Item(102, 8).Value = 24308.945
$ws.Cells.Item(102, 9).Value = 7971.6875
$ws.Cells.Item(102, 10).Value = 155007
$ws.Cells.Item(102, 11).Value = 7971.6875
$ws.Cells.Item(102, 12).Value = 155007
$ws.Cells.Item(102, 13).Value = -6349.6875
$ws.Cells.Item(102, 14).Value = -158251

$ws.Cells.Item(107, 8).Value = 721.4286
$ws.Cells.Item(107, 9).Value = 798.5714
$ws.Cells.Item(107, 10).Value = 644.2857
$ws.Cells.Item(107, 11).Value = 798.5714
$ws.Cells.Item(107, 12).Value = 644.2857
$ws.Cells.Item(107, 13).Value = 1121.4286
$ws.Cells.Item(107, 14).Value = -4484.2857

$ws.Cells.Item(113, 8).Value = 828.5714
$ws.Cells.Item(113, 9).Value = 780
$ws.Cells.Item(113, 10).Value = 950
$ws.Cells.Item(113, 11).Value = 780
$ws.Cells.Item(113, 12).Value = 950
$ws.Cells.Item(113, 13).Value = 1390
$ws.Cells.Item(113, 14).Value = -5290

$ws.Cells.Item(122, 8).Value = 2815.9092
$ws.Cells.Item(122, 9).Value = 2152.1428
$ws.Cells.Item(122, 10).Value = 3305
$ws.Cells.Item(122, 11).Value = 6456.428400000001
$ws.Cells.Item(122, 12).Value = 9915
$ws.Cells.Item(122, 13).Value = -4006.428400000001
$ws.Cells.Item(122, 14).Value = -14815

$ws.Cells.Item(126, 8).Value = 1990.909
$ws.Cells.Item(126, 9).Value = 1240
$ws.Cells.Item(126, 10).Value = 2157.7778
$ws.Cells.Item(126, 11).Value = 3720
$ws.Cells.Item(126, 12).Value = 6473.3334
$ws.Cells.Item(126, 13).Value = -1250
$ws.Cells.Item(126, 14).Value = -11413.3334

$ws.Cells.Item(132, 8).Value = 50539.43
$ws.Cells.Item(132, 9).Value = 2187.7144
$ws.Cells.Item(132, 10).Value = 74715.28999999999
$ws.Cells.Item(132, 11).Value = 6563.1432
$ws.Cells.Item(132, 12).Value = 224145.87
$ws.Cells.Item(132, 13).Value = -4033.1432
$ws.Cells.Item(132, 14).Value = -229205.87

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3830.5293
$ws.Cells.Item(7, 9).Value = 4409
$ws.Cells.Item(7, 10).Value = 2770
$ws.Cells.Item(7, 11).Value = 4409
$ws.Cells.Item(7, 12).Value = 2770
$ws.Cells.Item(7, 13).Value = -4297
$ws.Cells.Item(7, 14).Value = -2994

$ws.Cells.Item(40, 8).Value = 74879.21000000001
$ws.Cells.Item(40, 9).Value = 2340.5715
$ws.Cells.Item(40, 11).Value = 2340.5715
$ws.Cells.Item(40, 13).Value = -2204.5715

$ws.Cells.Item(55, 8).Value = 626.5294
$ws.Cells.Item(55, 9).Value = 448.0909
$ws.Cells.Item(55, 11).Value = 448.0909
$ws.Cells.Item(55, 13).Value = -275.0909

$ws.Cells.Item(61, 8).Value = 3116.15
$ws.Cells.Item(61, 9).Value = 2585.4285
$ws.Cells.Item(61, 10).Value = 3401.923
$ws.Cells.Item(61, 11).Value = 2585.4285
$ws.Cells.Item(61, 12).Value = 3401.923
$ws.Cells.Item(61, 13).Value = -2383.4285
$ws.Cells.Item(61, 14).Value = -3805.923

$ws.Cells.Item(93, 8).Value = 1434.6
$ws.Cells.Item(93, 9).Value = 851
$ws.Cells.Item(93, 10).Value = 2310
$ws.Cells.Item(93, 11).Value = 851
$ws.Cells.Item(93, 12).Value = 2310
$ws.Cells.Item(93, 13).Value = 397
$ws.Cells.Item(93, 14).Value = -4806

$ws.Cells.Item(113, 8).Value = 3116.15
$ws.Cells.Item(113, 9).Value = 2585.4285
$ws.Cells.Item(113, 10).Value = 3401.923
$ws.Cells.Item(113, 11).Value = 2585.4285
$ws.Cells.Item(113, 12).Value = 3401.923
$ws.Cells.Item(113, 13).Value = -415.4285
$ws.Cells.Item(113, 14).Value = -7741.923

$ws.Cells.Item(122, 8).Value = 2298.195
$ws.Cells.Item(122, 9).Value = 2018.8276
$ws.Cells.Item(122, 10).Value = 2973.3333
$ws.Cells.Item(122, 11).Value = 6056.4828
$ws.Cells.Item(122, 12).Value = 8919.999899999999
$ws.Cells.Item(122, 13).Value = -3606.4828
$ws.Cells.Item(122, 14).Value = -13819.9999

$ws.Cells.Item(126, 8).Value = 3830.5293
$ws.Cells.Item(126, 9).Value = 4409
$ws.Cells.Item(126, 10).Value = 2770
$ws.Cells.Item(126, 11).Value = 13227
$ws.Cells.Item(126, 12).Value = 8310
$ws.Cells.Item(126, 13).Value = -10757
$ws.Cells.Item(126, 14).Value = -13250

$ws.Cells.Item(136, 8).Value = 528697.0600000001
$ws.Cells.Item(136, 9).Value = 1113189.4
$ws.Cells.Item(136, 10).Value = 2654
$ws.Cells.Item(136, 11).Value = 3339568.2
$ws.Cells.Item(136, 12).Value = 7962
$ws.Cells.Item(136, 13).Value = -3337018.2
$ws.Cells.Item(136, 14).Value = -13062

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 5010.9375
$ws.Cells.Item(132, 9).Value = 1556.8334
$ws.Cells.Item(132, 10).Value = 15373.25
$ws.Cells.Item(132, 11).Value = 4670.5002
$ws.Cells.Item(132, 12).Value = 46119.75
$ws.Cells.Item(132, 13).Value = -2140.5002
$ws.Cells.Item(132, 14).Value = -51179.75
